$d = $word.ActiveDocument

# The document's single section has a primary header/footer and a
# first-page header/footer (both "exist"). Each one carries one inline
# picture:
#   - Headers -> the orange BTec logo, currently named "image2.jpg"
#     -> should become "image1.jpg"
#   - Footers -> the Pearson logo, currently named "image1.png"
#     -> should become "image2.png"
# Walk every section/header/footer combination and rename whichever
# inline picture is present.

foreach ($sec in $d.Sections) {

    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $hdr = $sec.Headers($hi)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ftr = $sec.Footers($fi)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
